$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Round row 5 (B5:AH5) values to 2 decimal places ("custom accuracy")
$ws.Cells.Item(5,2).Value  = 15.28    # B5
$ws.Cells.Item(5,3).Value  = 11.3     # C5
$ws.Cells.Item(5,4).Value  = 1.06     # D5
$ws.Cells.Item(5,5).Value  = 33.44    # E5
$ws.Cells.Item(5,6).Value  = 27.01    # F5
$ws.Cells.Item(5,7).Value  = 11.96    # G5
$ws.Cells.Item(5,8).Value  = 46.12    # H5
$ws.Cells.Item(5,9).Value  = 18.62    # I5
$ws.Cells.Item(5,10).Value = 8.19     # J5
$ws.Cells.Item(5,11).Value = 12.05    # K5
$ws.Cells.Item(5,12).Value = 13.4     # L5
$ws.Cells.Item(5,13).Value = 14.28    # M5
$ws.Cells.Item(5,14).Value = 3.86     # N5
$ws.Cells.Item(5,15).Value = 12.03    # O5
$ws.Cells.Item(5,16).Value = 17.05    # P5
$ws.Cells.Item(5,17).Value = 10.26    # Q5
$ws.Cells.Item(5,18).Value = 0.75     # R5
$ws.Cells.Item(5,19).Value = 0.69     # S5
$ws.Cells.Item(5,20).Value = 175.79   # T5
$ws.Cells.Item(5,21).Value = 33.62    # U5
$ws.Cells.Item(5,22).Value = 11.1     # V5
$ws.Cells.Item(5,23).Value = 22.48    # W5
$ws.Cells.Item(5,24).Value = 11.78    # X5
$ws.Cells.Item(5,25).Value = 1.88     # Y5
$ws.Cells.Item(5,26).Value = 22.54    # Z5
$ws.Cells.Item(5,27).Value = 9.81     # AA5
$ws.Cells.Item(5,28).Value = 8.75     # AB5
$ws.Cells.Item(5,29).Value = 10.29    # AC5
$ws.Cells.Item(5,30).Value = 14.04    # AD5
$ws.Cells.Item(5,31).Value = 0.54     # AE5
$ws.Cells.Item(5,32).Value = 41.97    # AF5
$ws.Cells.Item(5,33).Value = 6.2      # AG5
$ws.Cells.Item(5,34).Value = 13.88    # AH5

# 2) Remove row 6 (the old last data row) entirely - 1000-point dataset trim
$ws.Rows.Item(6).Delete()

# 3) Narrow a handful of data columns (C, K, Q, V, X) from width 8 to width 7,
#    matching the width already used by neighboring columns such as D.
$narrowWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Columns.Item(3).ColumnWidth  = $narrowWidth
$ws.Columns.Item(11).ColumnWidth = $narrowWidth
$ws.Columns.Item(17).ColumnWidth = $narrowWidth
$ws.Columns.Item(22).ColumnWidth = $narrowWidth
$ws.Columns.Item(24).ColumnWidth = $narrowWidth
